$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4 and 5 (F:H) need the same fill/border formatting already present on
# row 2 / row 3's F:H cells (style ids 9/11/11) before we populate them.
$ws.Range("F2:H2").Copy()
$ws.Range("F4:H5").PasteSpecial(-4122)

# New "Login" API method row (row 3)
$ws.Range("F3").Value = "Login "
# New "Logout" API method row (row 4)
$ws.Range("F4").Value = "Logout"
# New "Chek logged in" API method row (row 5)
$ws.Range("F5").Value = "Chek logged in "

# URLs (set in this order to reproduce the shared-string append order)
$ws.Range("G5").Value = "http://localhost:3005/users/loggedIn"
$ws.Range("G4").Value = "http://localhost:3005/users/logout"
$ws.Range("G3").Value = "http://localhost:3005/users/login"

# HTTP methods
$ws.Range("H3").Value = "POST"
$ws.Range("H4").Value = "GET"
$ws.Range("H5").Value = "GET"

# Update the saved view/selection
$ws.Range("H12").Select()
